# Icons.xlsx — add the "Xuat - Nhap - Ton kho" (inventory In/Out/Stock)
# ribbon icons to the documentation table, rows 18-30.
#
# Columns: A=MODULE  B=GUI  C=Control  D=Control(detail)  E=TEXT/CAPTION  F=FILE NAME
#
# NOTE: cell writes below are ordered so that brand-new strings are first
# introduced to the workbook in the same left-to-right / top-to-bottom
# sequence they end up in; this keeps the generated shared-string table
# stable and minimal (each unique string is interned exactly once, in
# first-use order), matching how the original data was authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: Xuat - Nhap - Ton kho (ribbon page)
$ws.Range("A18").Value = "VnsErp2025"
$ws.Range("B18").Value = "FormMain"
$ws.Range("C18").Value = "XuatNhapTonKhoRibbonPage"
$ws.Range("D18").Value = "XuatNhapTonKhoRibbonPage"
$ws.Range("E18").Value = "Xuất - Nhập - Tồn kho"
$ws.Range("F18").Value = "inventory.png"

# Row 19: Nhap kho group - Nhap bao hanh
$ws.Range("A19").Value = "VnsErp2025"
$ws.Range("B19").Value = "FormMain"
$ws.Range("D19").Value = "NhapBaoHanhBarButtonItem"
$ws.Range("C19").Value = "NhapKhoRibbonPageGroup"
$ws.Range("E19").Value = "Nhập bảo hành"
$ws.Range("F19").Value = "insurance.svg"

# Row 20: Nhap kho group - Nhap hang thuong mai
$ws.Range("A20").Value = "VnsErp2025"
$ws.Range("B20").Value = "FormMain"
$ws.Range("C20").Value = "NhapKhoRibbonPageGroup"
$ws.Range("F20").Value = "replacement.svg"
$ws.Range("E20").Value = "Nhập hàng thương mại"
$ws.Range("D20").Value = "NhapHangThuongMaiBarButtonItem"

# Row 21: Nhap kho group - Nhap luu chuyen kho
$ws.Range("A21").Value = "VnsErp2025"
$ws.Range("B21").Value = "FormMain"
$ws.Range("C21").Value = "NhapKhoRibbonPageGroup"
$ws.Range("D21").Value = "NhapLuuChuyenKhoBarButtonItem"
$ws.Range("E21").Value = "Nhập lưu chuyển kho"
$ws.Range("F21").Value = "supplier.svg"

# Row 22: Nhap kho group - Nhap noi bo
$ws.Range("A22").Value = "VnsErp2025"
$ws.Range("B22").Value = "FormMain"
$ws.Range("C22").Value = "NhapKhoRibbonPageGroup"
$ws.Range("D22").Value = "NhapNoiBoBarButtonItem"
$ws.Range("E22").Value = "Nhập nội bộ"
$ws.Range("F22").Value = "inventory (1).svg"

# Row 23: Nhap kho group - Nhap thiet bi muon - thue
$ws.Range("A23").Value = "VnsErp2025"
$ws.Range("B23").Value = "FormMain"
$ws.Range("C23").Value = "NhapKhoRibbonPageGroup"
$ws.Range("D23").Value = "NhapThietBiMuonBarButtonItem"
$ws.Range("E23").Value = "Nhập thiết bị mượn - thuê"
$ws.Range("F23").Value = "hardware.svg"

# Row 24: Xuat kho group - Xuat hang bao hanh
$ws.Range("A24").Value = "VnsErp2025"
$ws.Range("B24").Value = "FormMain"
$ws.Range("D24").Value = "XuatBaoHanhBarButtonItem"
$ws.Range("E24").Value = "Xuất hàng bảo hành"
$ws.Range("F24").Value = "insurance.svg"

# Row 25: Xuat kho group - Xuat hang thuong mai
$ws.Range("A25").Value = "VnsErp2025"
$ws.Range("B25").Value = "FormMain"
$ws.Range("D25").Value = "XuatHangThuongMaiBarButtonItem"
$ws.Range("E25").Value = "Xuất hàng thương mại"
$ws.Range("F25").Value = "replacement.svg"

# Row 26: Xuat kho group - Xuat luu chuyen kho
$ws.Range("A26").Value = "VnsErp2025"
$ws.Range("B26").Value = "FormMain"
$ws.Range("D26").Value = "XuatLuuChuyenKhoBarButtonItem"
$ws.Range("E26").Value = "Xuất lưu chuyển kho"
$ws.Range("F26").Value = "supplier.svg"

# Row 27: Xuat kho group - Xuat noi bo
$ws.Range("A27").Value = "VnsErp2025"
$ws.Range("B27").Value = "FormMain"
$ws.Range("D27").Value = "XuatNoiBoBarButtonItem"
$ws.Range("E27").Value = "Xuất nội bộ"
$ws.Range("F27").Value = "inventory (1).svg"

# Row 28: Xuat kho group - Xuat thiet bi muon - thue
$ws.Range("A28").Value = "VnsErp2025"
$ws.Range("B28").Value = "FormMain"
$ws.Range("D28").Value = "XuatChoThueMuonBarButtonItem"
$ws.Range("E28").Value = "Xuất thiết bị mượn - thuê"
$ws.Range("F28").Value = "hardware.svg"

# Back-fill the "XuatKhoRibbonPageGroup" control name shared by rows 24-28
$ws.Range("C24").Value = "XuatKhoRibbonPageGroup"
$ws.Range("C25").Value = "XuatKhoRibbonPageGroup"
$ws.Range("C26").Value = "XuatKhoRibbonPageGroup"
$ws.Range("C27").Value = "XuatKhoRibbonPageGroup"
$ws.Range("C28").Value = "XuatKhoRibbonPageGroup"

# Row 29: Stock In/Out history group - Phieu xuat kho
$ws.Range("D29").Value = "StockInOutMasterHistoryBarButtonItem"
$ws.Range("C29").Value = "StockInOutRibbonPageGroup"
$ws.Range("E29").Value = "Phiếu xuất kho"
$ws.Range("F29").Value = "invoice.svg"

# Row 30: Stock In/Out history group - San pham - dich vu
$ws.Range("D30").Value = "StockInOutProductHistoryBarButtonItem"
$ws.Range("E30").Value = "Sản phẩm - dịch vụ"
$ws.Range("F30").Value = "product-management.svg"

# Restore the author's last cell selection on the sheet
$ws.Range("K28").Select() | Out-Null
